$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 1181, shifting the existing
# 1181-1221 block (and everything after) down to 1183-1223.
$ws.Range("A1181:A1182").EntireRow.Insert()

# New row 1181: Acelga / Primera, week of 2022-06-10 (44722)
$ws.Cells.Item(1181, 1).Value2 = 6
$ws.Cells.Item(1181, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1181, 3).Value2 = "Metropolitana"
$ws.Cells.Item(1181, 4).Value2 = 44722
$ws.Cells.Item(1181, 5).Value2 = 13
$ws.Cells.Item(1181, 6).Value2 = 100112009
$ws.Cells.Item(1181, 7).Value2 = "Acelga"
$ws.Cells.Item(1181, 8).Value2 = "Sin especificar"
$ws.Cells.Item(1181, 9).Value2 = "Primera"
$ws.Cells.Item(1181, 10).Value2 = 150
$ws.Cells.Item(1181, 11).Value2 = 16000
$ws.Cells.Item(1181, 12).Value2 = 16000
$ws.Cells.Item(1181, 13).Value2 = 16000
$ws.Cells.Item(1181, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(1181, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(1181, 16).Value2 = 5333
$ws.Cells.Item(1181, 17).Value2 = 3
$ws.Cells.Item(1181, 18).Value2 = "Hortaliza"

# New row 1182: Acelga / Segunda, week of 2022-06-10 (44722)
$ws.Cells.Item(1182, 1).Value2 = 6
$ws.Cells.Item(1182, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1182, 3).Value2 = "Metropolitana"
$ws.Cells.Item(1182, 4).Value2 = 44722
$ws.Cells.Item(1182, 5).Value2 = 13
$ws.Cells.Item(1182, 6).Value2 = 100112009
$ws.Cells.Item(1182, 7).Value2 = "Acelga"
$ws.Cells.Item(1182, 8).Value2 = "Sin especificar"
$ws.Cells.Item(1182, 9).Value2 = "Segunda"
$ws.Cells.Item(1182, 10).Value2 = 90
$ws.Cells.Item(1182, 11).Value2 = 12000
$ws.Cells.Item(1182, 12).Value2 = 12000
$ws.Cells.Item(1182, 13).Value2 = 12000
$ws.Cells.Item(1182, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(1182, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(1182, 16).Value2 = 4000
$ws.Cells.Item(1182, 17).Value2 = 3
$ws.Cells.Item(1182, 18).Value2 = "Hortaliza"
